$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D header (row 1) ---
$ws.Range("D1").Value = "Manufacturer"
$ws.Range("D1").HorizontalAlignment = -4131   # xlHAlignLeft -> matches existing style (s="1")

# --- New manufacturer note on the Fan row (row 3) ---
$ws.Range("D3").Value = "Arctic"

# --- New row 6: Jar / 1 / 1590ml / Weck ---
$ws.Range("A6").Value = "Jar"
$ws.Range("A6").HorizontalAlignment = -4131

$ws.Range("B6").Value = 1
$ws.Range("B6").HorizontalAlignment = -4131

$ws.Range("C6").Value = "1590ml"
$ws.Range("C6").HorizontalAlignment = -4131

$ws.Range("D6").Value = "Weck"

# --- Column D width ---
$ws.Columns.Item(4).ColumnWidth = 15.29

# --- Move the active selection (matches the saved cursor position in the diff) ---
[void]$ws.Range("N24").Select()
